$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "test"
$ws.Range("C19").Value = "\Testdata\Non_Oncology\DataFiles\Protocol_Page\DownloadProtocol\Download_Protocol_Data.xlsx"
$ws.Range("B19").Value = "download_protocol_excel"

$ws.Range("B19").Select()
